# Natmi following Dr Hou advice
# Recompute the LR-pair stats for Jag2-Notch2 (per-edge ligand/receptor expressing-cell
# counts and derived expression/specificity weights) on the active sheet, rows 2-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.97098566666667
$ws.Range("H2").Value = 38.91295700000001
$ws.Range("I2").Value = 0.7291028508134716
$ws.Range("J2").Value = 0.7291028508134717
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 27.53580066666666
$ws.Range("N2").Value = 82.60740199999999
$ws.Range("O2").Value = 0.2054887285464767
$ws.Range("P2").Value = 0.2054887285464768
$ws.Range("Q2").Value = 357.1664757675238
$ws.Range("R2").Value = 3214.498281907714
$ws.Range("S2").Value = 0.1498224177932718
$ws.Range("T2").Value = 0.1498224177932718
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.97098566666667
$ws.Range("H3").Value = 38.91295700000001
$ws.Range("I3").Value = 0.7291028508134716
$ws.Range("J3").Value = 0.7291028508134717
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.50472933333333
$ws.Range("N3").Value = 82.514188
$ws.Range("O3").Value = 0.2052568555438283
$ws.Range("P3").Value = 0.2052568555438283
$ws.Range("Q3").Value = 356.763449948213
$ws.Range("R3").Value = 3210.871049533916
$ws.Range("S3").Value = 0.1496533585260141
$ws.Range("T3").Value = 0.1496533585260142
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.97098566666667
$ws.Range("H4").Value = 38.91295700000001
$ws.Range("I4").Value = 0.7291028508134716
$ws.Range("J4").Value = 0.7291028508134717
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 45.12975566666668
$ws.Range("N4").Value = 135.389267
$ws.Range("O4").Value = 0.3367854170582615
$ws.Range("P4").Value = 0.3367854170582616
$ws.Range("Q4").Value = 585.3774138925023
$ws.Range("R4").Value = 5268.39672503252
$ws.Range("S4").Value = 0.2455512076895825
$ws.Range("T4").Value = 0.2455512076895826
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.97098566666667
$ws.Range("H5").Value = 38.91295700000001
$ws.Range("I5").Value = 0.7291028508134716
$ws.Range("J5").Value = 0.7291028508134717
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.831228
$ws.Range("N5").Value = 101.493684
$ws.Range("O5").Value = 0.2524689988514334
$ws.Range("P5").Value = 0.2524689988514334
$ws.Range("Q5").Value = 438.8243734737321
$ws.Range("R5").Value = 3949.419361263589
$ws.Range("S5").Value = 0.1840758668046032
$ws.Range("T5").Value = 0.1840758668046032
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.047813333333333
$ws.Range("H6").Value = 3.14344
$ws.Range("I6").Value = 0.05889789011308234
$ws.Range("J6").Value = 0.05889789011308236
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.53580066666666
$ws.Range("N6").Value = 82.60740199999999
$ws.Range("O6").Value = 0.2054887285464767
$ws.Range("P6").Value = 0.2054887285464768
$ws.Range("Q6").Value = 28.85237908254222
$ws.Range("R6").Value = 259.67141174288
$ws.Range("S6").Value = 0.01210285255340739
$ws.Range("T6").Value = 0.0121028525534074
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.047813333333333
$ws.Range("H7").Value = 3.14344
$ws.Range("I7").Value = 0.05889789011308234
$ws.Range("J7").Value = 0.05889789011308236
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 27.50472933333333
$ws.Range("N7").Value = 82.514188
$ws.Range("O7").Value = 0.2052568555438283
$ws.Range("P7").Value = 0.2052568555438283
$ws.Range("Q7").Value = 28.81982212519111
$ws.Range("R7").Value = 259.37839912672
$ws.Range("S7").Value = 0.01208919572277721
$ws.Range("T7").Value = 0.01208919572277722
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.047813333333333
$ws.Range("H8").Value = 3.14344
$ws.Range("I8").Value = 0.05889789011308234
$ws.Range("J8").Value = 0.05889789011308236
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 45.12975566666668
$ws.Range("N8").Value = 135.389267
$ws.Range("O8").Value = 0.3367854170582615
$ws.Range("P8").Value = 0.3367854170582616
$ws.Range("Q8").Value = 47.28755971760889
$ws.Range("R8").Value = 425.5880374584801
$ws.Range("S8").Value = 0.0198359504855861
$ws.Range("T8").Value = 0.0198359504855861
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.047813333333333
$ws.Range("H9").Value = 3.14344
$ws.Range("I9").Value = 0.05889789011308234
$ws.Range("J9").Value = 0.05889789011308236
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.831228
$ws.Range("N9").Value = 101.493684
$ws.Range("O9").Value = 0.2524689988514334
$ws.Range("P9").Value = 0.2524689988514334
$ws.Range("Q9").Value = 35.44881178144
$ws.Range("R9").Value = 319.03930603296
$ws.Range("S9").Value = 0.01486989135131164
$ws.Range("T9").Value = 0.01486989135131164
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.312552
$ws.Range("H10").Value = 3.937656
$ws.Range("I10").Value = 0.07377892703252469
$ws.Range("J10").Value = 0.0737789270325247
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 27.53580066666666
$ws.Range("N10").Value = 82.60740199999999
$ws.Range("O10").Value = 0.2054887285464767
$ws.Range("P10").Value = 0.2054887285464768
$ws.Range("Q10").Value = 36.14217023663466
$ws.Range("R10").Value = 325.279532129712
$ws.Range("S10").Value = 0.01516073790943678
$ws.Range("T10").Value = 0.01516073790943679
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.312552
$ws.Range("H11").Value = 3.937656
$ws.Range("I11").Value = 0.07377892703252469
$ws.Range("J11").Value = 0.0737789270325247
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 27.50472933333333
$ws.Range("N11").Value = 82.514188
$ws.Range("O11").Value = 0.2052568555438283
$ws.Range("P11").Value = 0.2052568555438283
$ws.Range("Q11").Value = 36.10138749592533
$ws.Range("R11").Value = 324.912487463328
$ws.Range("S11").Value = 0.01514363056809357
$ws.Range("T11").Value = 0.01514363056809357
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.312552
$ws.Range("H12").Value = 3.937656
$ws.Range("I12").Value = 0.07377892703252469
$ws.Range("J12").Value = 0.0737789270325247
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 45.12975566666668
$ws.Range("N12").Value = 135.389267
$ws.Range("O12").Value = 0.3367854170582615
$ws.Range("P12").Value = 0.3367854170582616
$ws.Range("Q12").Value = 59.23515105979467
$ws.Range("R12").Value = 533.1163595381521
$ws.Range("S12").Value = 0.02484766671075987
$ws.Range("T12").Value = 0.02484766671075988
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.312552
$ws.Range("H13").Value = 3.937656
$ws.Range("I13").Value = 0.07377892703252469
$ws.Range("J13").Value = 0.0737789270325247
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.831228
$ws.Range("N13").Value = 101.493684
$ws.Range("O13").Value = 0.2524689988514334
$ws.Range("P13").Value = 0.2524689988514334
$ws.Range("Q13").Value = 44.405245973856
$ws.Range("R13").Value = 399.647213764704
$ws.Range("S13").Value = 0.01862689184423446
$ws.Range("T13").Value = 0.01862689184423446
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.458986333333333
$ws.Range("H14").Value = 7.376958999999999
$ws.Range("I14").Value = 0.1382203320409214
$ws.Range("J14").Value = 0.1382203320409214
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 27.53580066666666
$ws.Range("N14").Value = 82.60740199999999
$ws.Range("O14").Value = 0.2054887285464767
$ws.Range("P14").Value = 0.2054887285464768
$ws.Range("Q14").Value = 67.71015751672421
$ws.Range("R14").Value = 609.391417650518
$ws.Range("S14").Value = 0.02840272029036077
$ws.Range("T14").Value = 0.02840272029036078
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.458986333333333
$ws.Range("H15").Value = 7.376958999999999
$ws.Range("I15").Value = 0.1382203320409214
$ws.Range("J15").Value = 0.1382203320409214
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 27.50472933333333
$ws.Range("N15").Value = 82.514188
$ws.Range("O15").Value = 0.2052568555438283
$ws.Range("P15").Value = 0.2052568555438283
$ws.Range("Q15").Value = 67.6337535326991
$ws.Range("R15").Value = 608.703781794292
$ws.Range("S15").Value = 0.02837067072694337
$ws.Range("T15").Value = 0.02837067072694338
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.458986333333333
$ws.Range("H16").Value = 7.376958999999999
$ws.Range("I16").Value = 0.1382203320409214
$ws.Range("J16").Value = 0.1382203320409214
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 45.12975566666668
$ws.Range("N16").Value = 135.389267
$ws.Range("O16").Value = 0.3367854170582615
$ws.Range("P16").Value = 0.3367854170582616
$ws.Range("Q16").Value = 110.9734524110059
$ws.Range("R16").Value = 998.7610716990531
$ws.Range("S16").Value = 0.04655059217233309
$ws.Range("T16").Value = 0.04655059217233311
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.458986333333333
$ws.Range("H17").Value = 7.376958999999999
$ws.Range("I17").Value = 0.1382203320409214
$ws.Range("J17").Value = 0.1382203320409214
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.831228
$ws.Range("N17").Value = 101.493684
$ws.Range("O17").Value = 0.2524689988514334
$ws.Range("P17").Value = 0.2524689988514334
$ws.Range("Q17").Value = 83.19052729188401
$ws.Range("R17").Value = 748.714745626956
$ws.Range("S17").Value = 0.03489634885128411
$ws.Range("T17").Value = 0.03489634885128412
